$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 42
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 15
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 8
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 12
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 80
$ws.Range("X2").Value = 53
